$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 583; everything from the old row 583
# downward shifts down by one (old row 583 -> new row 584, ..., old
# row 627 -> new row 628).
$ws.Rows.Item(583).Insert()

# Populate the newly inserted row 583 with the new weekly price record.
$ws.Range("A583").Value = 8
$ws.Range("B583").Value = "Terminal La Palmera de La Serena"
$ws.Range("C583").Value = "Coquimbo"
$ws.Range("D583").Value = 45265
$ws.Range("E583").Value = 4
$ws.Range("F583").Value = 100112032
$ws.Range("G583").Value = "Zapallo italiano"
$ws.Range("H583").Value = "Sin especificar"
$ws.Range("I583").Value = "Primera"
$ws.Range("J583").Value = 500
$ws.Range("K583").Value = 8000
$ws.Range("L583").Value = 10000
$ws.Range("M583").Value = 9000
$ws.Range("N583").Value = "`$/caja 60 unidades"
$ws.Range("O583").Value = "Provincia de Limarí"
$ws.Range("P583").Value = 150
$ws.Range("Q583").Value = 60
$ws.Range("R583").Value = "Hortaliza"
